# "Moins bourrin, ajout du tableau, conso RTC"
# All modules but the RTC no longer draw any measurable current
# (consoOn/consoOff/tension reset to 0); the RTC row gets its real
# measured consumption and is highlighted as "Good"/"Satisfaisant".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,5,6,7,8,9 (ARDUINO, SIGFOX, GPS, DS18B20, DHT11, LCD, HX711):
# consoOn / consoOff / tension all go to 0.
foreach ($r in @(2, 3, 5, 6, 7, 8, 9)) {
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}

# Row 4 (RTC) keeps a real tension (5) and gets its measured consumption.
$ws.Cells.Item(4, 2).Value = 0.0015
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 5

# Highlight the RTC row with the built-in "Good" ("Satisfaisant") cell style.
$ws.Range("A4:D4").Style = "Good"

# Re-apply the scientific number format on the consumption columns so the
# highlighted row keeps showing values the same way as the rest of the table.
$ws.Cells.Item(4, 2).NumberFormat = "0.00E+00"
$ws.Cells.Item(4, 3).NumberFormat = "0.00E+00"

# Move the active selection to C5.
$ws.Range("C5").Select()
